# Scheduled market-data refresh: updates currentAveragePrice* / Leve cost-profit
# columns (H:N) on each job sheet with freshly pulled Universalis prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 622.75
$ws.Range("J32").Value = 630.3333
$ws.Range("L32").Value = 630.3333
$ws.Range("N32").Value = -1282.3333

$ws.Range("H92").Value = 1198.4762
$ws.Range("I92").Value = 1062.8235
$ws.Range("J92").Value = 1775
$ws.Range("K92").Value = 1062.8235
$ws.Range("L92").Value = 1775
$ws.Range("M92").Value = 185.1765
$ws.Range("N92").Value = -4271

$ws.Range("H99").Value = 1311.5714
$ws.Range("I99").Value = 1233.3334
$ws.Range("K99").Value = 3700.0002
$ws.Range("M99").Value = -2202.0002

$ws.Range("H100").Value = 100002500
$ws.Range("I100").Value = 200000000
$ws.Range("K100").Value = 200000000
$ws.Range("M100").Value = -199999459

$ws.Range("H113").Value = 3712.25
$ws.Range("J113").Value = 4309.8
$ws.Range("L113").Value = 4309.8
$ws.Range("N113").Value = -10817.8

$ws.Range("H115").Value = 1546.52
$ws.Range("I115").Value = 814.4545000000001
$ws.Range("K115").Value = 2443.3635
$ws.Range("M115").Value = -876.3635000000004

$ws.Range("H138").Value = 2876.16
$ws.Range("I138").Value = 780.4474
$ws.Range("J138").Value = 4160.629
$ws.Range("K138").Value = 2341.3422
$ws.Range("L138").Value = 12481.887
$ws.Range("M138").Value = 2798.6578
$ws.Range("N138").Value = -22761.887

$ws.Range("H141").Value = 4302.567
$ws.Range("I141").Value = 4297.3965
$ws.Range("K141").Value = 12892.1895
$ws.Range("M141").Value = -7712.1895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6325.4463
$ws.Range("I32").Value = 5095.4873
$ws.Range("J32").Value = 9147.117
$ws.Range("K32").Value = 5095.4873
$ws.Range("L32").Value = 9147.117
$ws.Range("M32").Value = -4808.4873
$ws.Range("N32").Value = -9721.117

$ws.Range("H63").Value = 12595228
$ws.Range("I63").Value = 15392835
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 15392835
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -15392149
$ws.Range("N63").Value = -7372

$ws.Range("H66").Value = 12595228
$ws.Range("I66").Value = 15392835
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 76964175
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -76960743
$ws.Range("N66").Value = -36864

$ws.Range("H74").Value = 2485.52
$ws.Range("I74").Value = 2317.372
$ws.Range("K74").Value = 2317.372
$ws.Range("M74").Value = -1443.372

$ws.Range("H77").Value = 2485.52
$ws.Range("I77").Value = 2317.372
$ws.Range("K77").Value = 11586.86
$ws.Range("M77").Value = -7218.859999999999

$ws.Range("H97").Value = 735.37036
$ws.Range("I97").Value = 717.5
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 717.5
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -221.5
$ws.Range("N97").Value = -2192

$ws.Range("H102").Value = 1588.7727
$ws.Range("I102").Value = 1514.8823
$ws.Range("K102").Value = 1514.8823
$ws.Range("M102").Value = 107.1177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 22316.666
$ws.Range("J60").Value = 22316.666
$ws.Range("L60").Value = 22316.666
$ws.Range("N60").Value = -23514.666

$ws.Range("H94").Value = 913.25
$ws.Range("I94").Value = 985.9
$ws.Range("J94").Value = 550
$ws.Range("K94").Value = 985.9
$ws.Range("L94").Value = 550
$ws.Range("M94").Value = -534.9
$ws.Range("N94").Value = -1452

$ws.Range("H99").Value = 4233.3335
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4233.3335
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 4233.3335
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -7229.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35719068
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 35719068
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 35719068
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -35719658

$ws.Range("H34").Value = 35719068
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 35719068
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 35719068
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -35719472

$ws.Range("H105").Value = 2035.5264
$ws.Range("I105").Value = 1837.0714
$ws.Range("J105").Value = 2591.2
$ws.Range("K105").Value = 1837.0714
$ws.Range("L105").Value = 2591.2
$ws.Range("M105").Value = -90.07140000000004
$ws.Range("N105").Value = -6085.2

$ws.Range("H122").Value = 3527.8
$ws.Range("I122").Value = 1592.6666
$ws.Range("J122").Value = 4357.143
$ws.Range("K122").Value = 4777.9998
$ws.Range("L122").Value = 13071.429
$ws.Range("M122").Value = -2327.9998
$ws.Range("N122").Value = -17971.429

$ws.Range("H135").Value = 37390
$ws.Range("J135").Value = 37390
$ws.Range("L135").Value = 37390
$ws.Range("N135").Value = -47530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 521.48
$ws.Range("I113").Value = 505.7857
$ws.Range("J113").Value = 541.4545000000001
$ws.Range("K113").Value = 1517.3571
$ws.Range("L113").Value = 1624.3635
$ws.Range("M113").Value = 652.6428999999998
$ws.Range("N113").Value = -5964.3635

$ws.Range("H136").Value = 3697.8
$ws.Range("J136").Value = 3403.5557
$ws.Range("L136").Value = 10210.6671
$ws.Range("N136").Value = -20410.6671

$ws.Range("H137").Value = 2655.7273
$ws.Range("J137").Value = 4038.923
$ws.Range("L137").Value = 12116.769
$ws.Range("N137").Value = -22316.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1105
$ws.Range("I97").Value = 1105
$ws.Range("K97").Value = 1105
$ws.Range("M97").Value = -609

$ws.Range("H102").Value = 2373.9546
$ws.Range("I102").Value = 1716.6
$ws.Range("J102").Value = 3782.5715
$ws.Range("K102").Value = 1716.6
$ws.Range("L102").Value = 3782.5715
$ws.Range("M102").Value = -94.59999999999991
$ws.Range("N102").Value = -7026.5715

$ws.Range("H123").Value = 12723.444
$ws.Range("J123").Value = 12723.444
$ws.Range("L123").Value = 12723.444
$ws.Range("N123").Value = -17623.444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 58800
$ws.Range("J36").Value = 58800
$ws.Range("L36").Value = 58800
$ws.Range("N36").Value = -59924

$ws.Range("H93").Value = 7409591
$ws.Range("I93").Value = 15874742
$ws.Range("J93").Value = 2583.25
$ws.Range("K93").Value = 15874742
$ws.Range("L93").Value = 2583.25
$ws.Range("M93").Value = -15873494
$ws.Range("N93").Value = -5079.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 722.3333
$ws.Range("I100").Value = 722.3333
$ws.Range("K100").Value = 1444.6666
$ws.Range("M100").Value = -903.6666

$ws.Range("H136").Value = 3139.125
$ws.Range("I136").Value = 755.4091
$ws.Range("J136").Value = 8383.299999999999
$ws.Range("K136").Value = 2266.2273
$ws.Range("L136").Value = 25149.9
$ws.Range("M136").Value = 283.7727
$ws.Range("N136").Value = -30249.9
